$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.52749240398407
$ws.Range("B1").Value = 1.740251183509827
$ws.Range("C1").Value = 2.163393020629883
$ws.Range("D1").Value = 2.217680931091309
$ws.Range("E1").Value = 1.367837071418762
